# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is numeric-looking text (must be forced to Text
#     format first so Excel does not silently convert them to numbers and
#     lose the exact printed representation, e.g. trailing zeros). ---
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.70'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4901'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.92'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2433'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05988'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06730'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.67'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.431'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5873'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.40'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.48'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006370'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.941'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.119'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.829'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '135.41'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.836'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.450'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.45'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '100.08'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08097'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.714'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.384'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04345'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9997'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.660'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.019'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6024'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.740'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.043'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '102.37'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01481'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7954'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3810'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.116'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.011'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05094'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.06'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1050'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.27'

# --- Remaining cells (coin names, links, price strings that Excel cannot
#     misinterpret as numbers, and the volume/percentage strings). ---
$ws.Range("D2").Value = '25.685.88'
$ws.Range("E2").Value = '  -3.53%  '
$ws.Range("D3").Value = '1.739.03'
$ws.Range("E3").Value = '  -5.64%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  -8.01%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -7.01%  '
$ws.Range("E8").Value = '  -7.04%  '
$ws.Range("E9").Value = '  -22.71%  '
$ws.Range("E10").Value = '  -12.06%  '
$ws.Range("D11").Value = '1.738.19'
$ws.Range("E11").Value = '  -5.61%  '
$ws.Range("E12").Value = '  -13.25%  '
$ws.Range("E13").Value = '  -21.67%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E14").Value = '  -11.57%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("E15").Value = '  -24.87%  '
$ws.Range("E16").Value = '  -13.21%  '
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '25.719.64'
$ws.Range("E19").Value = '  -3.48%  '
$ws.Range("E20").Value = '  -17.16%  '
$ws.Range("E21").Value = '  -19.70%  '
$ws.Range("D22").Value = '1.954.57'
$ws.Range("E22").Value = '  -5.75%  '
$ws.Range("E24").Value = '  -14.38%  '
$ws.Range("E25").Value = '  -15.99%  '
$ws.Range("E26").Value = '  -5.02%  '
$ws.Range("E27").Value = '  -16.77%  '
$ws.Range("E28").Value = '  -13.71%  '
$ws.Range("E29").Value = '  -14.92%  '
$ws.Range("E30").Value = '  -9.75%  '
$ws.Range("E31").Value = '  -7.26%  '
$ws.Range("E32").Value = '  -11.26%  '
$ws.Range("E33").Value = '  -16.97%  '
$ws.Range("E34").Value = '  -10.91%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("E36").Value = '  -6.98%  '
$ws.Range("E37").Value = '  -10.65%  '
$ws.Range("E38").Value = '  -17.60%  '
$ws.Range("E39").Value = '  -11.35%  '
$ws.Range("E40").Value = '  -9.51%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -6.95%  '
$ws.Range("E43").Value = '  -14.42%  '
$ws.Range("E44").Value = '  -11.08%  '
$ws.Range("E45").Value = '  -20.53%  '
$ws.Range("E46").Value = '  -13.66%  '
$ws.Range("E47").Value = '  -21.60%  '
$ws.Range("E48").Value = '  -12.42%  '
$ws.Range("E49").Value = '  -13.51%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E50").Value = '  -15.35%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E51").Value = '  -12.59%  '
